$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "M1"
$ws.Range("B2").Value = "Ltf"
$ws.Range("C2").Value = "Lrp11"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.245485
$ws.Range("H2").Value = 0.736455
$ws.Range("I2").Value = 0.01511172246591349
$ws.Range("J2").Value = 0.01511172246591349
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5915956666666666
$ws.Range("N2").Value = 1.774787
$ws.Range("O2").Value = 0.127834779324208
$ws.Range("P2").Value = 0.127834779324208
$ws.Range("Q2").Value = 0.1452278622316666
$ws.Range("R2").Value = 1.307050760085
$ws.Range("S2").Value = 0.001931803706638728
$ws.Range("T2").Value = 0.001931803706638728

# Row 3
$ws.Range("A3").Value = "M1"
$ws.Range("B3").Value = "Ltf"
$ws.Range("C3").Value = "Lrp11"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.245485
$ws.Range("H3").Value = 0.736455
$ws.Range("I3").Value = 0.01511172246591349
$ws.Range("J3").Value = 0.01511172246591349
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.070248333333333
$ws.Range("N3").Value = 3.210745
$ws.Range("O3").Value = 0.2312643030072366
$ws.Range("P3").Value = 0.2312643030072365
$ws.Range("Q3").Value = 0.2627299121083334
$ws.Range("R3").Value = 2.364569208975
$ws.Range("S3").Value = 0.003494801963318282
$ws.Range("T3").Value = 0.003494801963318281

# Row 4
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Ltf"
$ws.Range("C4").Value = "Lrp11"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.245485
$ws.Range("H4").Value = 0.736455
$ws.Range("I4").Value = 0.01511172246591349
$ws.Range("J4").Value = 0.01511172246591349
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.603607
$ws.Range("N4").Value = 1.810821
$ws.Range("O4").Value = 0.1304302448297411
$ws.Range("P4").Value = 0.1304302448297411
$ws.Range("Q4").Value = 0.148176464395
$ws.Range("R4").Value = 1.333588179555
$ws.Range("S4").Value = 0.001971025661028196
$ws.Range("T4").Value = 0.001971025661028195

# Row 5
$ws.Range("A5").Value = "M1"
$ws.Range("B5").Value = "Ltf"
$ws.Range("C5").Value = "Lrp11"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.245485
$ws.Range("H5").Value = 0.736455
$ws.Range("I5").Value = 0.01511172246591349
$ws.Range("J5").Value = 0.01511172246591349
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.7085566666666666
$ws.Range("N5").Value = 2.12567
$ws.Range("O5").Value = 0.1531082633386932
$ws.Range("P5").Value = 0.1531082633386932
$ws.Range("Q5").Value = 0.1739400333166666
$ws.Range("R5").Value = 1.56546029985
$ws.Range("S5").Value = 0.002313729582812329
$ws.Range("T5").Value = 0.002313729582812329

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Ltf"
$ws.Range("C6").Value = "Lrp11"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.245485
$ws.Range("H6").Value = 0.736455
$ws.Range("I6").Value = 0.01511172246591349
$ws.Range("J6").Value = 0.01511172246591349
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.447354
$ws.Range("N6").Value = 1.342062
$ws.Range("O6").Value = 0.0966663603065637
$ws.Range("P6").Value = 0.09666636030656368
$ws.Range("Q6").Value = 0.10981869669
$ws.Range("R6").Value = 0.9883682702099998
$ws.Range("S6").Value = 0.001460795208742787
$ws.Range("T6").Value = 0.001460795208742786

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Ltf"
$ws.Range("C7").Value = "Lrp11"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.245485
$ws.Range("H7").Value = 0.736455
$ws.Range("I7").Value = 0.01511172246591349
$ws.Range("J7").Value = 0.01511172246591349
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.206453
$ws.Range("N7").Value = 3.619359
$ws.Range("O7").Value = 0.2606960491935574
$ws.Range("P7").Value = 0.2606960491935574
$ws.Range("Q7").Value = 0.296166114705
$ws.Range("R7").Value = 2.665495032345
$ws.Range("S7").Value = 0.003939566343373171
$ws.Range("T7").Value = 0.00393956634337317

# Row 8
$ws.Range("A8").Value = "Neutro"
$ws.Range("B8").Value = "Ltf"
$ws.Range("C8").Value = "Lrp11"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 15.99918866666667
$ws.Range("H8").Value = 47.99756600000001
$ws.Range("I8").Value = 0.9848882775340866
$ws.Range("J8").Value = 0.9848882775340865
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5915956666666666
$ws.Range("N8").Value = 1.774787
$ws.Range("O8").Value = 0.127834779324208
$ws.Range("P8").Value = 0.127834779324208
$ws.Range("Q8").Value = 9.465050685382446
$ws.Range("R8").Value = 85.185456168442
$ws.Range("S8").Value = 0.1259029756175693
$ws.Range("T8").Value = 0.1259029756175693

# Row 9
$ws.Range("A9").Value = "Neutro"
$ws.Range("B9").Value = "Ltf"
$ws.Range("C9").Value = "Lrp11"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 15.99918866666667
$ws.Range("H9").Value = 47.99756600000001
$ws.Range("I9").Value = 0.9848882775340866
$ws.Range("J9").Value = 0.9848882775340865
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.070248333333333
$ws.Range("N9").Value = 3.210745
$ws.Range("O9").Value = 0.2312643030072366
$ws.Range("P9").Value = 0.2312643030072365
$ws.Range("Q9").Value = 17.12310500518556
$ws.Range("R9").Value = 154.10794504667
$ws.Range("S9").Value = 0.2277695010439183
$ws.Range("T9").Value = 0.2277695010439182

# Row 10
$ws.Range("A10").Value = "Neutro"
$ws.Range("B10").Value = "Ltf"
$ws.Range("C10").Value = "Lrp11"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 15.99918866666667
$ws.Range("H10").Value = 47.99756600000001
$ws.Range("I10").Value = 0.9848882775340866
$ws.Range("J10").Value = 0.9848882775340865
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.603607
$ws.Range("N10").Value = 1.810821
$ws.Range("O10").Value = 0.1304302448297411
$ws.Range("P10").Value = 0.1304302448297411
$ws.Range("Q10").Value = 9.657222273520668
$ws.Range("R10").Value = 86.915000461686
$ws.Range("S10").Value = 0.1284592191687129
$ws.Range("T10").Value = 0.1284592191687129

# Row 11
$ws.Range("A11").Value = "Neutro"
$ws.Range("B11").Value = "Ltf"
$ws.Range("C11").Value = "Lrp11"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 15.99918866666667
$ws.Range("H11").Value = 47.99756600000001
$ws.Range("I11").Value = 0.9848882775340866
$ws.Range("J11").Value = 0.9848882775340865
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.7085566666666666
$ws.Range("N11").Value = 2.12567
$ws.Range("O11").Value = 0.1531082633386932
$ws.Range("P11").Value = 0.1531082633386932
$ws.Range("Q11").Value = 11.33633179102445
$ws.Range("R11").Value = 102.02698611922
$ws.Range("S11").Value = 0.1507945337558809
$ws.Range("T11").Value = 0.1507945337558808

# Row 12
$ws.Range("A12").Value = "Neutro"
$ws.Range("B12").Value = "Ltf"
$ws.Range("C12").Value = "Lrp11"
$ws.Range("D12").Value = "Neutro"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 15.99918866666667
$ws.Range("H12").Value = 47.99756600000001
$ws.Range("I12").Value = 0.9848882775340866
$ws.Range("J12").Value = 0.9848882775340865
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.447354
$ws.Range("N12").Value = 1.342062
$ws.Range("O12").Value = 0.0966663603065637
$ws.Range("P12").Value = 0.09666636030656368
$ws.Range("Q12").Value = 7.157301046788
$ws.Range("R12").Value = 64.415709421092
$ws.Range("S12").Value = 0.09520556509782092
$ws.Range("T12").Value = 0.0952055650978209

# Row 13
$ws.Range("A13").Value = "Neutro"
$ws.Range("B13").Value = "Ltf"
$ws.Range("C13").Value = "Lrp11"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 15.99918866666667
$ws.Range("H13").Value = 47.99756600000001
$ws.Range("I13").Value = 0.9848882775340866
$ws.Range("J13").Value = 0.9848882775340865
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.206453
$ws.Range("N13").Value = 3.619359
$ws.Range("O13").Value = 0.2606960491935574
$ws.Range("P13").Value = 0.2606960491935574
$ws.Range("Q13").Value = 19.302269164466
$ws.Range("R13").Value = 173.720422480194
$ws.Range("S13").Value = 0.2567564828501843
$ws.Range("T13").Value = 0.2567564828501843
